# Weekly update: add this week's "Haba" price rows (Mercado Mayorista Lo
# Valledor de Santiago) at the top of the data block, pushing the previously
# most-recent rows down by 4 rows.
#
# Before: rows 114-140 held the data for 2021-09-29 (44468) and earlier.
# After : 4 brand-new rows (114-117) for 2021-10-07 (44476) are inserted
#         right after the existing header/data row 113, and everything that
#         used to be rows 114-140 is now 118-144 (values unchanged).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 114..140 down by inserting 4 fresh rows at 114.
$ws.Range("A114:A117").EntireRow.Insert()

# Common values shared by every row in this data block.
$mercado = "Mercado Mayorista Lo Valledor de Santiago"
$region  = "Metropolitana"
$codreg  = 13
$catId   = 100112026
$categoria = "Haba"
$variedad  = "Sin especificar"
$unidad    = "`$/malla 25 kilos"
$kgUnidades = 25
$clasificacion = "Hortaliza"
$fecha = 44476

$newRows = @(
    @{ Row=114; Calidad="Primera"; Volumen=800; PMin=6000;  PMax=7000;  PProm=6438; Origen="Región Metropolitana";  PrecioKg=258 },
    @{ Row=115; Calidad="Primera"; Volumen=400; PMin=6000;  PMax=7000;  PProm=6425; Origen="Región de Coquimbo";    PrecioKg=257 },
    @{ Row=116; Calidad="Segunda"; Volumen=250; PMin=5000;  PMax=5000;  PProm=5000; Origen="Región Metropolitana";  PrecioKg=200 },
    @{ Row=117; Calidad="Segunda"; Volumen=120; PMin=5000;  PMax=5000;  PProm=5000; Origen="Región de Coquimbo";    PrecioKg=200 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value  = 6
    $ws.Cells.Item($row, 2).Value  = $mercado
    $ws.Cells.Item($row, 3).Value  = $region
    $ws.Cells.Item($row, 4).Value  = $fecha
    $ws.Cells.Item($row, 5).Value  = $codreg
    $ws.Cells.Item($row, 6).Value  = $catId
    $ws.Cells.Item($row, 7).Value  = $categoria
    $ws.Cells.Item($row, 8).Value  = $variedad
    $ws.Cells.Item($row, 9).Value  = $r.Calidad
    $ws.Cells.Item($row, 10).Value = $r.Volumen
    $ws.Cells.Item($row, 11).Value = $r.PMin
    $ws.Cells.Item($row, 12).Value = $r.PMax
    $ws.Cells.Item($row, 13).Value = $r.PProm
    $ws.Cells.Item($row, 14).Value = $unidad
    $ws.Cells.Item($row, 15).Value = $r.Origen
    $ws.Cells.Item($row, 16).Value = $r.PrecioKg
    $ws.Cells.Item($row, 17).Value = $kgUnidades
    $ws.Cells.Item($row, 18).Value = $clasificacion
}
